# Add 2022-Q4 data: insert a new quarter sheet and update the summary sheet.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Update the "总计" (totals) sheet: prepend a new 2022-Q4 row and shift the
#    existing quarters down by one.
# ---------------------------------------------------------------------------
$total = $wb.Worksheets.Item(1)

$totalData = @(
    @(0, "2022-Q4", 7, 1.04),
    @(1, "2022-Q3", 8, 0.93),
    @(2, "2022-Q2", 21, 0.6899999999999999),
    @(3, "2022-Q1", 4, 0.43),
    @(4, "2021-Q4", 3, 0.43),
    @(5, "2021-Q3", 1, 0.09),
    @(6, "2021-Q2", 1, 0.1)
)

$r = 2
foreach ($row in $totalData) {
    $total.Cells.Item($r, 1).Value = $row[0]
    $total.Cells.Item($r, 2).Value = $row[1]
    $total.Cells.Item($r, 3).Value = $row[2]
    $total.Cells.Item($r, 4).Value = $row[3]
    $r = $r + 1
}

# Row 8 (2021-Q2) is brand new — give column A the same look (bold, centered,
# boxed) as the rest of the index column above it.
$total.Range("A7").Copy()
$total.Range("A8").PasteSpecial(-4122)
$total.Range("A8").Value = 6

# ---------------------------------------------------------------------------
# 2. Insert a brand-new "2022-Q4" sheet right after "总计" (i.e. before the
#    sheet that is currently "2022-Q3") with the quarter's fund holdings.
# ---------------------------------------------------------------------------
$beforeSheet = $wb.Worksheets.Item(2)
$q4 = $wb.Worksheets.Add($beforeSheet)
$q4.Name = "2022-Q4"

# NOTE: worksheet handles captured before an Add() can resolve to the wrong
# sheet once indices shift (the old "2022-Q3" sheet moved from position 2 to
# 3) -- look it up again by name now that the insert has happened.
$oldQ3 = $wb.Worksheets.Item("2022-Q3")

# Match the header/index-column look of the neighbouring quarter sheet.
$oldQ3.Range("B1:H1").Copy()
$q4.Range("B1:H1").PasteSpecial(-4122)
$oldQ3.Range("A2:A8").Copy()
$q4.Range("A2:A8").PasteSpecial(-4122)

$q4.Range("B1").Value = "基金代码"
$q4.Range("C1").Value = "基金名称"
$q4.Range("D1").Value = "基金规模"
$q4.Range("E1").Value = "股票总仓位"
$q4.Range("F1").Value = "仓位占比"
$q4.Range("G1").Value = "持有市值(亿元)"
$q4.Range("H1").Value = "仓位排名"

# Columns B-G hold numeric-looking text (fund codes, names and percentages
# stored verbatim) — force text so things like leading zeros survive.
$q4.Range("B2:G8").NumberFormat = "@"

$q4Data = @(
    @(0, "519700", "交银主题优选混合A",       "26.24", "73.61", "2.17", "0.5694", 8),
    @(1, "013884", "交银主题优选混合C",       "11.57", "73.61", "2.17", "0.2511", 8),
    @(2, "001628", "招商体育文化休闲股票A",    "2.33",  "93.03", "4.93", "0.1149", 6),
    @(3, "159855", "银华中证影视主题ETF",      "1.01",  "97.80", "4.20", "0.0424", 7),
    @(4, "516620", "国泰中证影视主题ETF",      "0.71",  "98.01", "3.90", "0.0277", 9),
    @(5, "003397", "银华体育文化灵活配置混合", "0.53",  "87.74", "2.87", "0.0152", 10),
    @(6, "015395", "招商体育文化休闲股票C",    "0.29",  "93.03", "4.93", "0.0143", 6)
)

$r = 2
foreach ($row in $q4Data) {
    $q4.Cells.Item($r, 1).Value = $row[0]
    $q4.Cells.Item($r, 2).Value = $row[1]
    $q4.Cells.Item($r, 3).Value = $row[2]
    $q4.Cells.Item($r, 4).Value = $row[3]
    $q4.Cells.Item($r, 5).Value = $row[4]
    $q4.Cells.Item($r, 6).Value = $row[5]
    $q4.Cells.Item($r, 7).Value = $row[6]
    $q4.Cells.Item($r, 8).Value = $row[7]
    $r = $r + 1
}

# ---------------------------------------------------------------------------
# 3. Restore the originally-active tab (last sheet, "2021-Q2") since adding a
#    sheet shifts Excel's focus to the newly created one.
# ---------------------------------------------------------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$lastSheet.Activate()
$wb.Worksheets.Item("总计").Activate()
